# Insert a new weekly data row for "Acelga" / Terminal La Palmera de La Serena
# at row 664, shifting all existing rows (664-741) down to (665-742).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 664 (pushes old 664..741 down to 665..742)
$ws.Rows.Item(664).EntireRow.Insert()

# Populate the newly inserted row 664 with the new record's data.
$ws.Cells.Item(664, 1).Value = 8
$ws.Cells.Item(664, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(664, 3).Value = "Coquimbo"
$ws.Cells.Item(664, 4).Value = 45212
$ws.Cells.Item(664, 5).Value = 4
$ws.Cells.Item(664, 6).Value = 100112009
$ws.Cells.Item(664, 7).Value = "Acelga"
$ws.Cells.Item(664, 8).Value = "Sin especificar"
$ws.Cells.Item(664, 9).Value = "Primera"
$ws.Cells.Item(664, 10).Value = 1500
$ws.Cells.Item(664, 11).Value = 600
$ws.Cells.Item(664, 12).Value = 700
$ws.Cells.Item(664, 13).Value = 650
$ws.Cells.Item(664, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(664, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(664, 16).Value = 325
$ws.Cells.Item(664, 17).Value = 2
$ws.Cells.Item(664, 18).Value = "Hortaliza"

Write-Host "Row inserted and populated."
